# Insert a new weekly price record as row 203, pushing the existing
# rows 203-318 down to 204-319 (dimension grows from A1:R318 to A1:R319).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(203).Insert()

$ws.Range("A203").Value = 9
$ws.Range("B203").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C203").Value = "Metropolitana"
$ws.Range("D203").Value = 44813
$ws.Range("E203").Value = 13
$ws.Range("F203").Value = 300000001
$ws.Range("G203").Value = "Rabanito"
$ws.Range("H203").Value = "Sin especificar"
$ws.Range("I203").Value = "Primera"
$ws.Range("J203").Value = 7000
$ws.Range("K203").Value = 2500
$ws.Range("L203").Value = 3000
$ws.Range("M203").Value = 2750
$ws.Range("N203").Value = "`$/cien unidades (volumen en unidades)"
$ws.Range("O203").Value = "Provincia de Chacabuco"
$ws.Range("P203").Value = 28
$ws.Range("Q203").Value = 100
$ws.Range("R203").Value = "Hortaliza"
